$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5708.607
$ws.Range("I32").Value = 5162.375
$ws.Range("J32").Value = 5927.1
$ws.Range("K32").Value = 5162.375
$ws.Range("L32").Value = 5927.1
$ws.Range("M32").Value = -4836.375
$ws.Range("N32").Value = -6579.1
$ws.Range("H45").Value = 6008.5
$ws.Range("J45").Value = 12000
$ws.Range("L45").Value = 36000
$ws.Range("N45").Value = -36384
$ws.Range("H103").Value = 2335.4285
$ws.Range("I103").Value = 1374
$ws.Range("J103").Value = 2720
$ws.Range("K103").Value = 4122
$ws.Range("L103").Value = 8160
$ws.Range("M103").Value = -3536
$ws.Range("N103").Value = -9332
$ws.Range("H112").Value = 1917.091
$ws.Range("J112").Value = 2008.8
$ws.Range("L112").Value = 6026.4
$ws.Range("N112").Value = -8242.4
$ws.Range("H116").Value = 8724.857
$ws.Range("I116").Value = 9727.546
$ws.Range("K116").Value = 9727.546
$ws.Range("M116").Value = -6285.546
$ws.Range("H138").Value = 3374.9512
$ws.Range("I138").Value = 1578.6316
$ws.Range("J138").Value = 3916.6985
$ws.Range("K138").Value = 4735.8948
$ws.Range("L138").Value = 11750.0955
$ws.Range("M138").Value = 404.1052
$ws.Range("N138").Value = -22030.0955
$ws.Range("H141").Value = 79576.55499999999
$ws.Range("I141").Value = 14523.625
$ws.Range("K141").Value = 43570.875
$ws.Range("M141").Value = -38390.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2886.2449
$ws.Range("I32").Value = 3988.0417
$ws.Range("J32").Value = 1828.52
$ws.Range("K32").Value = 3988.0417
$ws.Range("L32").Value = 1828.52
$ws.Range("M32").Value = -3701.0417
$ws.Range("N32").Value = -2402.52
$ws.Range("H88").Value = 2102.4333
$ws.Range("I88").Value = 2125.7727
$ws.Range("J88").Value = 2038.25
$ws.Range("K88").Value = 2125.7727
$ws.Range("L88").Value = 2038.25
$ws.Range("M88").Value = -1719.7727
$ws.Range("N88").Value = -2850.25
$ws.Range("H91").Value = 2102.4333
$ws.Range("I91").Value = 2125.7727
$ws.Range("J91").Value = 2038.25
$ws.Range("K91").Value = 2125.7727
$ws.Range("L91").Value = 2038.25
$ws.Range("M91").Value = -721.7727
$ws.Range("N91").Value = -4846.25
$ws.Range("H112").Value = 54500
$ws.Range("J112").Value = 54500
$ws.Range("L112").Value = 54500
$ws.Range("N112").Value = -57454
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3909.878
$ws.Range("I20").Value = 3472.2903
$ws.Range("J20").Value = 5266.4
$ws.Range("K20").Value = 3472.2903
$ws.Range("L20").Value = 5266.4
$ws.Range("M20").Value = -3225.2903
$ws.Range("N20").Value = -5760.4
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H86").Value = 9854.429
$ws.Range("I86").Value = 7003.4517
$ws.Range("J86").Value = 17889
$ws.Range("K86").Value = 7003.4517
$ws.Range("L86").Value = 17889
$ws.Range("M86").Value = -5880.4517
$ws.Range("N86").Value = -20135
$ws.Range("H89").Value = 9854.429
$ws.Range("I89").Value = 7003.4517
$ws.Range("J89").Value = 17889
$ws.Range("K89").Value = 35017.2585
$ws.Range("L89").Value = 89445
$ws.Range("M89").Value = -29401.2585
$ws.Range("N89").Value = -100677
$ws.Range("H99").Value = 4873.5
$ws.Range("I99").Value = 2498.3333
$ws.Range("K99").Value = 2498.3333
$ws.Range("M99").Value = -1000.3333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9804.406000000001
$ws.Range("I134").Value = 2309.32
$ws.Range("K134").Value = 6927.960000000001
$ws.Range("M134").Value = -4392.960000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 13281.944
$ws.Range("I56").Value = 13281.944
$ws.Range("K56").Value = 13281.944
$ws.Range("M56").Value = -12751.944
$ws.Range("H74").Value = 9000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 9000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H92").Value = 590.5714
$ws.Range("I92").Value = 245.71428
$ws.Range("J92").Value = 935.4286
$ws.Range("K92").Value = 737.14284
$ws.Range("L92").Value = 2806.2858
$ws.Range("M92").Value = 510.85716
$ws.Range("N92").Value = -5302.2858
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = ""
$ws.Range("H137").Value = 3958.0625
$ws.Range("I137").Value = 4348.385
$ws.Range("J137").Value = 2266.6667
$ws.Range("K137").Value = 13045.155
$ws.Range("L137").Value = 6800.000100000001
$ws.Range("M137").Value = -7945.155000000001
$ws.Range("N137").Value = -17000.0001
$ws.Range("H140").Value = 18829.143
$ws.Range("I140").Value = 18829.143
$ws.Range("K140").Value = 56487.429
$ws.Range("M140").Value = -51307.429
$ws.Range("H141").Value = 129464.75
$ws.Range("I141").Value = 5097.857
$ws.Range("K141").Value = 15293.571
$ws.Range("M141").Value = -10113.571
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 27781028
$ws.Range("I18").Value = 37039704
$ws.Range("K18").Value = 37039704
$ws.Range("M18").Value = -37039411
$ws.Range("H70").Value = 9852.352999999999
$ws.Range("I70").Value = 10045.923
$ws.Range("J70").Value = 9223.25
$ws.Range("K70").Value = 10045.923
$ws.Range("L70").Value = 9223.25
$ws.Range("M70").Value = -9775.923000000001
$ws.Range("N70").Value = -9763.25
$ws.Range("H73").Value = 9852.352999999999
$ws.Range("I73").Value = 10045.923
$ws.Range("J73").Value = 9223.25
$ws.Range("K73").Value = 10045.923
$ws.Range("L73").Value = 9223.25
$ws.Range("M73").Value = -9109.923000000001
$ws.Range("N73").Value = -11095.25
$ws.Range("H80").Value = 2473.5
$ws.Range("I80").Value = 2303
$ws.Range("J80").Value = 2700.8333
$ws.Range("K80").Value = 2303
$ws.Range("L80").Value = 2700.8333
$ws.Range("M80").Value = -1305
$ws.Range("N80").Value = -4696.8333
$ws.Range("H83").Value = 2473.5
$ws.Range("I83").Value = 2303
$ws.Range("J83").Value = 2700.8333
$ws.Range("K83").Value = 11515
$ws.Range("L83").Value = 13504.1665
$ws.Range("M83").Value = -6523
$ws.Range("N83").Value = -23488.1665
$ws.Range("H97").Value = 1618.9375
$ws.Range("I97").Value = 1747.2222
$ws.Range("K97").Value = 1747.2222
$ws.Range("M97").Value = -1251.2222
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = ""
$ws.Range("H119").Value = 29999
$ws.Range("J119").Value = 29999
$ws.Range("L119").Value = 29999
$ws.Range("N119").Value = -39675
$ws.Range("H122").Value = 1048.9524
$ws.Range("I122").Value = 925.4666999999999
$ws.Range("J122").Value = 1357.6666
$ws.Range("K122").Value = 2776.4001
$ws.Range("L122").Value = 4072.9998
$ws.Range("M122").Value = -326.4000999999998
$ws.Range("N122").Value = -8972.9998
$ws.Range("H132").Value = 220064.5
$ws.Range("I132").Value = 246452.56
$ws.Range("K132").Value = 739357.6799999999
$ws.Range("M132").Value = -736827.6799999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3999999.2
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3999999.2
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3999999.2
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -4000223.2
$ws.Range("H7").Value = 23905.285
$ws.Range("I7").Value = 30800.8
$ws.Range("J7").Value = 6666.5
$ws.Range("K7").Value = 30800.8
$ws.Range("L7").Value = 6666.5
$ws.Range("M7").Value = -30688.8
$ws.Range("N7").Value = -6890.5
$ws.Range("H22").Value = 5403.533
$ws.Range("I22").Value = 1942.8572
$ws.Range("K22").Value = 1942.8572
$ws.Range("M22").Value = -1647.8572
$ws.Range("H27").Value = 5403.533
$ws.Range("I27").Value = 1942.8572
$ws.Range("K27").Value = 1942.8572
$ws.Range("M27").Value = -1835.8572
$ws.Range("H43").Value = 24999.666
$ws.Range("I43").Value = 24999
$ws.Range("K43").Value = 24999
$ws.Range("M43").Value = -24806
$ws.Range("H69").Value = 100000
$ws.Range("I69").Value = 100000
$ws.Range("K69").Value = 100000
$ws.Range("M69").Value = -99189
$ws.Range("H72").Value = 100000
$ws.Range("I72").Value = 100000
$ws.Range("K72").Value = 300000
$ws.Range("M72").Value = -295944
$ws.Range("H100").Value = 43939.082
$ws.Range("I100").Value = 49858.953
$ws.Range("K100").Value = 49858.953
$ws.Range("M100").Value = -49317.953
$ws.Range("H126").Value = 23905.285
$ws.Range("I126").Value = 30800.8
$ws.Range("J126").Value = 6666.5
$ws.Range("K126").Value = 92402.39999999999
$ws.Range("L126").Value = 19999.5
$ws.Range("M126").Value = -89932.39999999999
$ws.Range("N126").Value = -24939.5
$ws.Range("H132").Value = 4013.0322
$ws.Range("I132").Value = 3362.238
$ws.Range("K132").Value = 10086.714
$ws.Range("M132").Value = -7556.714
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50000
$ws.Range("I2").Value = 50000
$ws.Range("K2").Value = 50000
$ws.Range("M2").Value = -49888
$ws.Range("H122").Value = 4354.6665
$ws.Range("I122").Value = 4024
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 12072
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -9622
$ws.Range("N122").Value = -25900
$ws.Range("H132").Value = 2338.68
$ws.Range("I132").Value = 2303.1304
$ws.Range("K132").Value = 6909.3912
$ws.Range("M132").Value = -4379.3912

Write-Host "Applied all edits"